$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Fix existing rows 16/17 (A16: 2 -> 3) ---
$ws.Range("A16").Value = 3

# --- Row 22: header/title row (copy of row 11) ---
$ws.Range("A22").Value = "st_level"
$ws.Range("B22").Value = "st_levelTable"
$ws.Range("C22").Value = '{"IsStringId":false,"IsGenItemClass":true,"JSONName":"st_levelJSON","IsGenEnum":true,"Path":"toanstt","IsSeparatedJSON":true}'

# --- Row 23: column header row (copy of row 12, but A is "Id" string) ---
$ws.Range("A23").Value = "Id"
$ws.Range("B23").Value = "float:row"
$ws.Range("C23").Value = "string:col"
$ws.Range("D23").Value = "is_boss"
$ws.Range("E23").Value = "myarray0"
$ws.Range("F23").Value = "myarray1"
$ws.Range("G23").Value = "myarray2"
$ws.Range("H23").Value = "bool:testfloat"
$ws.Range("I23").Value = "teststring"
$ws.Range("J23").Value = "array0"
$ws.Range("K23").Value = "array1"
$ws.Range("L23").Value = "TestEnum:testenum"
$ws.Range("M23").Value = "TestEnum:e0"
$ws.Range("N23").Value = "TestEnum:e1"
$ws.Range("O23").Value = "TestEnum:e2"

# --- Row 24: data row (copy of row 13, A becomes string id "r4r") ---
$ws.Range("A24").Value = "r4r"
$ws.Range("B24").Value = 4
$ws.Range("C24").Value = 4
$ws.Range("D24").Value = $true
$ws.Range("E24").Value = 23
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 2.5
$ws.Range("I24").Value = "asd"
$ws.Range("J24").Value = "string1"
$ws.Range("K24").Value = "strings2"
$ws.Range("L24").Value = "ENUM1:5"
$ws.Range("M24").Value = "ENUM3:67"
$ws.Range("N24").Value = "ENUM1"

# --- Row 25: data row (copy of row 14, A becomes string id "5r4") ---
$ws.Range("A25").Value = "5r4"
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = $false
$ws.Range("E25").Value = 2
$ws.Range("F25").Value = 3
$ws.Range("I25").Value = "asd"
$ws.Range("J25").Value = "string2"
$ws.Range("K25").Value = "strings3"
$ws.Range("L25").Value = "ENUM2"
$ws.Range("M25").Value = "ENUM2"
$ws.Range("N25").Value = "ENUM3"

# --- Row 26: data row (copy of row 15, A becomes string id "ffd") ---
$ws.Range("A26").Value = "ffd"
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 2
$ws.Range("F26").Value = 3
$ws.Range("I26").Value = "asd"
$ws.Range("J26").Value = "string3"
$ws.Range("K26").Value = "strings4"
$ws.Range("L26").Value = "ENUM2"
$ws.Range("M26").Value = "ENUM2"
$ws.Range("N26").Value = "ENUM2"

# --- Row 27: data row (copy of row 16, A becomes string id "45fe") ---
$ws.Range("A27").Value = "45fe"
$ws.Range("B27").Value = 4
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 2
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 2.4
$ws.Range("I27").Value = "df"
$ws.Range("J27").Value = "string4"
$ws.Range("K27").Value = "strings5"
$ws.Range("L27").Value = "ENUM1"
$ws.Range("M27").Value = "ENUM1"

# --- Row 28: data row (copy of row 17, A becomes string id "34f") ---
$ws.Range("A28").Value = "34f"
$ws.Range("B28").Value = 4
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 2
$ws.Range("H28").Value = 1.2
$ws.Range("I28").Value = "dsf"
$ws.Range("J28").Value = "string5"
$ws.Range("K28").Value = "strings6"
$ws.Range("L28").Value = "ENUM3"
$ws.Range("M28").Value = "ENUM3"

# --- Carry over the "quotePrefix" cell style used on column D in rows 3-7 (and 13-17) ---
$ws.Range("D3:D7").Copy()
$ws.Range("D24:D28").PasteSpecial(-4122)

# --- Update the active selection to match the new work location ---
$null = $ws.Range("F23").Select()
